$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.798.77"
$ws.Range("E2").Value = "  +3.61%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.802.74"
$ws.Range("E3").Value = "  +8.19%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "425.70"
$ws.Range("E5").Value = "  +8.16%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.42"
$ws.Range("E6").Value = "  +6.54%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.799.21"
$ws.Range("E7").Value = "  +8.40%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.609"
$ws.Range("E8").Value = "  +3.39%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.02%  "

$ws.Range("E10").Value = "  +7.07%  "

$ws.Range("E11").Value = "  +1.37%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000332"
$ws.Range("E12").Value = "  -1.45%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "41.40"
$ws.Range("E13").Value = "  +6.41%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.46"
$ws.Range("E14").Value = "  +13.34%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.437.50"
$ws.Range("E15").Value = "  +9.50%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.51"
$ws.Range("E16").Value = "  +21.63%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.863.24"
$ws.Range("E17").Value = "  +11.12%  "

$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.138"
$ws.Range("E18").Value = "  +1.38%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.94"
$ws.Range("E19").Value = "  +6.34%  "

$ws.Range("E20").Value = "  +7.93%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "66.154.28"
$ws.Range("E21").Value = "  +4.26%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "414.18"
$ws.Range("E22").Value = "  +4.65%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.98"
$ws.Range("E23").Value = "  +7.55%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.08"
$ws.Range("E24").Value = "  +4.56%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.11"
$ws.Range("E25").Value = "  +8.15%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "37.18"
$ws.Range("E26").Value = "  +10.69%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.13"
$ws.Range("E27").Value = "  +14.63%  "

$ws.Range("E28").Value = "  +9.91%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.47"
$ws.Range("E29").Value = "  +39.80%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.38"
$ws.Range("E30").Value = "  +2.69%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "14.15"
$ws.Range("E31").Value = "  +18.98%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "709.16"
$ws.Range("E32").Value = "  +5.42%  "

$ws.Range("E33").Value = "  +13.12%  "

$ws.Range("E34").Value = "  +6.04%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.10%  "

$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.74"
$ws.Range("E36").Value = "  +42.86%  "

$ws.Range("B37").Value = "InjectiveProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "38.63"
$ws.Range("E37").Value = "  +5.28%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.148"
$ws.Range("E38").Value = "  +0.39%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "55.59"
$ws.Range("E39").Value = "  +3.81%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0468"
$ws.Range("E40").Value = "  +7.12%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0725"
$ws.Range("E41").Value = "  +13.77%  "

$ws.Range("E42").Value = "  +4.22%  "

$ws.Range("E43").Value = "  +0.57%  "

$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.137"
$ws.Range("E44").Value = "  +4.51%  "

$ws.Range("B45").Value = "LidoDAOToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.38"
$ws.Range("E45").Value = "  +9.78%  "

$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.19"
$ws.Range("E46").Value = "  +3.63%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.319"
$ws.Range("E47").Value = "  +15.47%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.42"
$ws.Range("E48").Value = "  +44.30%  "

$ws.Range("E49").Value = "  +7.87%  "

$ws.Range("E50").Value = "  +5.03%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.83"
$ws.Range("E51").Value = "  +3.18%  "
